# Remove the "Myelinating Schwann cells" row from the Brain marker table.
# (cellName column is alphabetically sorted, so this is row 10 in the
# original sheet, right after "Microglial cells" and before
# "Neural Progenitor cells".)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(10).Delete()

# Leave the selection where the author's Excel session ended up.
$ws.Range("C21").Select()
